$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the rows that correspond to classes removed from the report
# (labels "32", "33", "42", "44" at original rows 17, 18, 20, 21).
# Delete from bottom to top so row indices of earlier rows stay valid.
$ws.Rows.Item(21).Delete()
$ws.Rows.Item(20).Delete()
$ws.Rows.Item(18).Delete()
$ws.Rows.Item(17).Delete()

# Update the remaining numeric cells (precision / recall / f1-score / support)
# with their new values, row by row.

# row 2 -> label "0"
$ws.Cells.Item(2, 2).Value = 0.8333333333333334
$ws.Cells.Item(2, 3).Value = 1
$ws.Cells.Item(2, 4).Value = 0.9090909090909091
$ws.Cells.Item(2, 5).Value = 5

# row 3 -> label "2"
$ws.Cells.Item(3, 2).Value = 1
$ws.Cells.Item(3, 3).Value = 0.2
$ws.Cells.Item(3, 4).Value = 0.3333333333333334
$ws.Cells.Item(3, 5).Value = 5

# row 4 -> label "5"
$ws.Cells.Item(4, 2).Value = 0.8571428571428571
$ws.Cells.Item(4, 3).Value = 1
$ws.Cells.Item(4, 4).Value = 0.923076923076923
$ws.Cells.Item(4, 5).Value = 6

# row 5 -> label "6"
$ws.Cells.Item(5, 2).Value = 1
$ws.Cells.Item(5, 3).Value = 1
$ws.Cells.Item(5, 4).Value = 1
$ws.Cells.Item(5, 5).Value = 3

# row 6 -> label "8"
$ws.Cells.Item(6, 2).Value = 1
$ws.Cells.Item(6, 3).Value = 0.8571428571428571
$ws.Cells.Item(6, 4).Value = 0.923076923076923
$ws.Cells.Item(6, 5).Value = 14

# row 7 -> label "11"
$ws.Cells.Item(7, 2).Value = 0.3333333333333333
$ws.Cells.Item(7, 3).Value = 0.5
$ws.Cells.Item(7, 4).Value = 0.4
$ws.Cells.Item(7, 5).Value = 2

# row 8 -> label "13"
$ws.Cells.Item(8, 2).Value = 1
$ws.Cells.Item(8, 3).Value = 1
$ws.Cells.Item(8, 4).Value = 1
$ws.Cells.Item(8, 5).Value = 3

# row 9 -> label "14"
$ws.Cells.Item(9, 2).Value = 1
$ws.Cells.Item(9, 3).Value = 0.6
$ws.Cells.Item(9, 4).Value = 0.7499999999999999
$ws.Cells.Item(9, 5).Value = 5

# row 10 -> label "15" (unchanged, but set explicitly for completeness)
$ws.Cells.Item(10, 2).Value = 1
$ws.Cells.Item(10, 3).Value = 0.6666666666666666
$ws.Cells.Item(10, 4).Value = 0.8
$ws.Cells.Item(10, 5).Value = 3

# row 11 -> label "16"
$ws.Cells.Item(11, 2).Value = 1
$ws.Cells.Item(11, 3).Value = 1
$ws.Cells.Item(11, 4).Value = 1
$ws.Cells.Item(11, 5).Value = 2

# row 12 -> label "18"
$ws.Cells.Item(12, 2).Value = 0.7142857142857143
$ws.Cells.Item(12, 3).Value = 0.8333333333333334
$ws.Cells.Item(12, 4).Value = 0.7692307692307692
$ws.Cells.Item(12, 5).Value = 6

# row 13 -> label "19" (unchanged, but set explicitly for completeness)
$ws.Cells.Item(13, 2).Value = 0.6666666666666666
$ws.Cells.Item(13, 3).Value = 1
$ws.Cells.Item(13, 4).Value = 0.8
$ws.Cells.Item(13, 5).Value = 2

# row 14 -> label "21"
$ws.Cells.Item(14, 2).Value = 0.75
$ws.Cells.Item(14, 3).Value = 1
$ws.Cells.Item(14, 4).Value = 0.8571428571428571
$ws.Cells.Item(14, 5).Value = 3

# row 15 -> label "27"
$ws.Cells.Item(15, 2).Value = 0.4
$ws.Cells.Item(15, 3).Value = 0.6666666666666666
$ws.Cells.Item(15, 4).Value = 0.5
$ws.Cells.Item(15, 5).Value = 3

# row 16 -> label "28"
$ws.Cells.Item(16, 2).Value = 0.6666666666666666
$ws.Cells.Item(16, 3).Value = 0.5
$ws.Cells.Item(16, 4).Value = 0.5714285714285715
$ws.Cells.Item(16, 5).Value = 4

# row 17 -> label "39" (was row 19 before the row deletions)
$ws.Cells.Item(17, 2).Value = 1
$ws.Cells.Item(17, 3).Value = 1
$ws.Cells.Item(17, 4).Value = 1
$ws.Cells.Item(17, 5).Value = 3

# row 18 -> label "46" (was row 22 before the row deletions)
$ws.Cells.Item(18, 2).Value = 0.75
$ws.Cells.Item(18, 3).Value = 1
$ws.Cells.Item(18, 4).Value = 0.8571428571428571
$ws.Cells.Item(18, 5).Value = 3

# row 19 -> label "50" (was row 23 before the row deletions)
$ws.Cells.Item(19, 2).Value = 0.8
$ws.Cells.Item(19, 3).Value = 1
$ws.Cells.Item(19, 4).Value = 0.888888888888889
$ws.Cells.Item(19, 5).Value = 4

# row 20 -> label "53" (was row 24 before the row deletions)
$ws.Cells.Item(20, 2).Value = 1
$ws.Cells.Item(20, 3).Value = 1
$ws.Cells.Item(20, 4).Value = 1
$ws.Cells.Item(20, 5).Value = 3

# row 21 -> label "accuracy" (was row 25 before the row deletions)
$ws.Cells.Item(21, 2).Value = 0.8227848101265823
$ws.Cells.Item(21, 3).Value = 0.8227848101265823
$ws.Cells.Item(21, 4).Value = 0.8227848101265823
$ws.Cells.Item(21, 5).Value = 0.8227848101265823

# row 22 -> label "macro avg" (was row 26 before the row deletions)
$ws.Cells.Item(22, 2).Value = 0.8300751879699249
$ws.Cells.Item(22, 3).Value = 0.8328320802005011
$ws.Cells.Item(22, 4).Value = 0.804337475390107
$ws.Cells.Item(22, 5).Value = 79

# row 23 -> label "weighted avg" (was row 27 before the row deletions)
$ws.Cells.Item(23, 2).Value = 0.8628089210367692
$ws.Cells.Item(23, 3).Value = 0.8227848101265823
$ws.Cells.Item(23, 4).Value = 0.8142173578882439
$ws.Cells.Item(23, 5).Value = 79
